# AutoCommit_23 апреля 2024 г. 10:11:52_SibNout2023
#
# Fills in missing ДЗ (homework) scores for three students who had
# zero/blank marks, matching the "Лаб_1" (column F) style/no-fill
# formatting that the rest of the sheet already uses once a score is
# entered. xlPasteFormats copies only the formatting (not the value) so
# a cell's score can be set afterwards without disturbing the copied
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-Score {
    param(
        [string]$Cell,
        [double]$Value
    )
    # Borrow column F's (unfilled/no colour) cell format for this cell,
    # then write the new score into it.
    $ws.Range("F4").Copy() | Out-Null
    $ws.Range($Cell).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($Cell).Value = $Value
}

# Борисов Никита (row 4): ДЗ_1 and ДЗ_3 were 0, Лаб_1 was blank.
Set-Score "C4" 5
Set-Score "E4" 5
$ws.Range("F4").Value = 5

# Карачун Анастасия (row 14): ДЗ_1..ДЗ_3 were 0, Лаб_1 was blank.
Set-Score "C14" 5
Set-Score "D14" 5
Set-Score "E14" 5
$ws.Range("F14").Value = 5

# Чебан Александра (row 27): ДЗ_3 was 0, Лаб_1 was blank.
Set-Score "E27" 5
$ws.Range("F27").Value = 5

$excel.CutCopyMode = $false

# Restore the view: scroll the frozen pane back up and park the
# selection on G5.
$ws.Range("G5").Select()
